# Optuna Attempt (go back with original)
# Updates forecast numbers on "Forecast Comparison" and the roll-up
# totals on "Summary" to the values produced by the new run.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: MyForecast (D), Inventory Coverage (H),
#     Seasonality Index (L) ---

$wsForecast.Range("D2").Value = 67
$wsForecast.Range("H2").Value = 3.98
$wsForecast.Range("L2").Value = 1.18

$wsForecast.Range("D3").Value = 67
$wsForecast.Range("H3").Value = 2.98
$wsForecast.Range("L3").Value = 0.95

$wsForecast.Range("H4").Value = 2.39
$wsForecast.Range("L4").Value = 1.12

$wsForecast.Range("H5").Value = 1.39
$wsForecast.Range("L5").Value = 0.85

$wsForecast.Range("D6").Value = 42
$wsForecast.Range("H6").Value = 0.52
$wsForecast.Range("L6").Value = 0.92

$wsForecast.Range("D7").Value = 50
$wsForecast.Range("L7").Value = 0.82

$wsForecast.Range("L8").Value = 1

$wsForecast.Range("L9").Value = 1.12

$wsForecast.Range("L10").Value = 1.11

$wsForecast.Range("L11").Value = 0.83

$wsForecast.Range("L12").Value = 1.18

$wsForecast.Range("L13").Value = 1.13

$wsForecast.Range("L14").Value = 1.16

$wsForecast.Range("L15").Value = 1.01

$wsForecast.Range("L16").Value = 1.13

$wsForecast.Range("L17").Value = 1.17

# --- Summary sheet: totals are stored as text (not numbers), so force
#     text formatting before writing, then clear the format override
#     back out to keep the default (General) cell style. ---

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $wsSummary.Range("B9")  "875"
Set-TextValue $wsSummary.Range("B10") "451"
Set-TextValue $wsSummary.Range("B11") "247"
Set-TextValue $wsSummary.Range("B12") "68"
Set-TextValue $wsSummary.Range("B14") "42"
